$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.618.80'
$ws.Range('E2').Value = '  +1.98%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.912.18'
$ws.Range('E3').Value = '  +3.91%  '
$ws.Range('E4').Value = '  +0.60%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.66'
$ws.Range('E5').Value = '  +5.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.634'
$ws.Range('E6').Value = '  +2.64%  '
$ws.Range('E7').Value = '  +0.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.34'
$ws.Range('E8').Value = '  +3.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.338'
$ws.Range('E9').Value = '  +3.40%  '
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.187.40'
$ws.Range('E12').Value = '  +3.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '12.49'
$ws.Range('E13').Value = '  +10.04%  '
$ws.Range('E14').Value = '  +3.58%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.882.64'
$ws.Range('E15').Value = '  +2.33%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.84'
$ws.Range('E16').Value = '  +4.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.576.73'
$ws.Range('E17').Value = '  +1.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '72.01'
$ws.Range('E18').Value = '  +3.14%  '
$ws.Range('E19').Value = '  +2.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '243.92'
$ws.Range('E20').Value = '  +1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.49'
$ws.Range('E21').Value = '  +2.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.94'
$ws.Range('E22').Value = '  +4.17%  '
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('E24').Value = '  +1.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.35'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.10'
$ws.Range('E26').Value = '  +25.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.48'
$ws.Range('E27').Value = '  +8.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.02'
$ws.Range('E28').Value = '  +3.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.125'
$ws.Range('E29').Value = '  +1.57%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.952'
$ws.Range('E30').Value = '  +27.29%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.10'
$ws.Range('E31').Value = '  +3.91%  '
$ws.Range('E32').Value = '  +2.75%  '
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.17'
$ws.Range('E34').Value = '  +6.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.72'
$ws.Range('E35').Value = '  +5.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.04'
$ws.Range('E36').Value = '  +4.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.33'
$ws.Range('E37').Value = '  +6.12%  '
$ws.Range('E38').Value = '  +5.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0205'
$ws.Range('E39').Value = '  +4.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '91.90'
$ws.Range('E40').Value = '  +2.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.360.94'
$ws.Range('E41').Value = '  +1.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.23'
$ws.Range('E42').Value = '  +4.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '49.19'
$ws.Range('E44').Value = '  +11.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.09'
$ws.Range('E45').Value = '  +19.25%  '
$ws.Range('E47').Value = '  +1.04%  '
$ws.Range('E48').Value = '  +1.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.67'
$ws.Range('E49').Value = '  +5.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.095.16'
$ws.Range('E50').Value = '  +3.27%  '
$ws.Range('E51').Value = '  +4.68%  '
